$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3849.4783
$ws.Range("I86").Value = 3950.8333
$ws.Range("J86").Value = 3813.7058
$ws.Range("K86").Value = 3950.8333
$ws.Range("L86").Value = 3813.7058
$ws.Range("M86").Value = -2827.8333
$ws.Range("N86").Value = -6059.7058
$ws.Range("H88").Value = 2473135.8
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2473135.8
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2473135.8
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -2473947.8
$ws.Range("H89").Value = 3849.4783
$ws.Range("I89").Value = 3950.8333
$ws.Range("J89").Value = 3813.7058
$ws.Range("K89").Value = 19754.1665
$ws.Range("L89").Value = 19068.529
$ws.Range("M89").Value = -14138.1665
$ws.Range("N89").Value = -30300.529
$ws.Range("H91").Value = 2473135.8
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2473135.8
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2473135.8
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -2475943.8
$ws.Range("H100").Value = 2675
$ws.Range("I100").Value = 2350
$ws.Range("K100").Value = 2350
$ws.Range("M100").Value = -1809
$ws.Range("H129").Value = 891.4231
$ws.Range("J129").Value = 961.6667
$ws.Range("L129").Value = 2885.0001
$ws.Range("N129").Value = -12885.0001
$ws.Range("H135").Value = 1045.2307
$ws.Range("I135").Value = 275.05
$ws.Range("K135").Value = 2475.45
$ws.Range("M135").Value = 59.54999999999973
$ws.Range("H138").Value = 1590.0344
$ws.Range("I138").Value = 1397.2307
$ws.Range("J138").Value = 1746.6875
$ws.Range("K138").Value = 4191.6921
$ws.Range("L138").Value = 5240.0625
$ws.Range("M138").Value = 948.3078999999998
$ws.Range("N138").Value = -15520.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3170.7576
$ws.Range("I32").Value = 2881.2458
$ws.Range("K32").Value = 2881.2458
$ws.Range("M32").Value = -2594.2458
$ws.Range("H61").Value = 938.8
$ws.Range("I61").Value = 938.8
$ws.Range("K61").Value = 938.8
$ws.Range("M61").Value = -726.8
$ws.Range("H74").Value = 822.3022999999999
$ws.Range("I74").Value = 475.17142
$ws.Range("K74").Value = 475.17142
$ws.Range("M74").Value = 398.82858
$ws.Range("H77").Value = 822.3022999999999
$ws.Range("I77").Value = 475.17142
$ws.Range("K77").Value = 2375.8571
$ws.Range("M77").Value = 1992.1429
$ws.Range("H110").Value = 1584.8
$ws.Range("I110").Value = 1065.6
$ws.Range("J110").Value = 2623.2
$ws.Range("K110").Value = 1065.6
$ws.Range("L110").Value = 2623.2
$ws.Range("M110").Value = 979.4000000000001
$ws.Range("N110").Value = -6713.2
$ws.Range("H122").Value = 1701.3846
$ws.Range("I122").Value = 1701.3846
$ws.Range("K122").Value = 5104.1538
$ws.Range("M122").Value = -2654.1538
$ws.Range("H132").Value = 1640.8
$ws.Range("I132").Value = 1327.9231
$ws.Range("K132").Value = 3983.7693
$ws.Range("M132").Value = -1453.7693
$ws.Range("H134").Value = 34144.11
$ws.Range("J134").Value = 34144.11
$ws.Range("L134").Value = 34144.11
$ws.Range("N134").Value = -44284.11
$ws.Range("H136").Value = 938.8
$ws.Range("I136").Value = 938.8
$ws.Range("K136").Value = 2816.4
$ws.Range("M136").Value = -266.3999999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1350.9615
$ws.Range("I31").Value = 1104.2778
$ws.Range("J31").Value = 1906
$ws.Range("K31").Value = 1104.2778
$ws.Range("L31").Value = 1906
$ws.Range("M31").Value = -809.2778000000001
$ws.Range("N31").Value = -2496
$ws.Range("H34").Value = 1350.9615
$ws.Range("I34").Value = 1104.2778
$ws.Range("J34").Value = 1906
$ws.Range("K34").Value = 1104.2778
$ws.Range("L34").Value = 1906
$ws.Range("M34").Value = -902.2778000000001
$ws.Range("N34").Value = -2310
$ws.Range("H58").Value = 1949.5294
$ws.Range("I58").Value = 1595.1666
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 1595.1666
$ws.Range("L58").Value = 2800
$ws.Range("M58").Value = -1392.1666
$ws.Range("N58").Value = -3206
$ws.Range("H134").Value = 1665.6061
$ws.Range("I134").Value = 1667.6666
$ws.Range("J134").Value = 1656.3334
$ws.Range("K134").Value = 5002.9998
$ws.Range("L134").Value = 4969.0002
$ws.Range("M134").Value = -2467.9998
$ws.Range("N134").Value = -10039.0002
$ws.Range("H136").Value = 1949.5294
$ws.Range("I136").Value = 1595.1666
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 4785.4998
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -2235.4998
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1727.9333
$ws.Range("I34").Value = 1325.8
$ws.Range("J34").Value = 1929
$ws.Range("K34").Value = 3977.4
$ws.Range("L34").Value = 5787
$ws.Range("M34").Value = -3893.4
$ws.Range("N34").Value = -5955
$ws.Range("H39").Value = 1899.4546
$ws.Range("J39").Value = 1599.3684
$ws.Range("L39").Value = 4798.1052
$ws.Range("N39").Value = -5386.1052
$ws.Range("H55").Value = 2801
$ws.Range("H113").Value = 605.1795
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 674.14813
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 2022.44439
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -6362.444390000001
$ws.Range("H129").Value = 37879496
$ws.Range("I129").Value = 41667136
$ws.Range("J129").Value = 27779124
$ws.Range("K129").Value = 125001408
$ws.Range("L129").Value = 83337372
$ws.Range("M129").Value = -124996408
$ws.Range("N129").Value = -83347372
$ws.Range("H131").Value = 13514701
$ws.Range("I131").Value = 250000270
$ws.Range("J131").Value = 1239.6
$ws.Range("K131").Value = 750000810
$ws.Range("L131").Value = 3718.8
$ws.Range("M131").Value = -749995770
$ws.Range("N131").Value = -13798.8
$ws.Range("H140").Value = 38051.93
$ws.Range("I140").Value = 54667.316
$ws.Range("J140").Value = 2975
$ws.Range("K140").Value = 164001.948
$ws.Range("L140").Value = 8925
$ws.Range("M140").Value = -158821.948
$ws.Range("N140").Value = -19285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 18645.154
$ws.Range("J46").Value = 18645.154
$ws.Range("L46").Value = 18645.154
$ws.Range("N46").Value = -18957.154
$ws.Range("H70").Value = 56253004
$ws.Range("I70").Value = 125002000
$ws.Range("J70").Value = 33336672
$ws.Range("K70").Value = 125002000
$ws.Range("L70").Value = 33336672
$ws.Range("M70").Value = -125001730
$ws.Range("N70").Value = -33337212
$ws.Range("H73").Value = 56253004
$ws.Range("I73").Value = 125002000
$ws.Range("J73").Value = 33336672
$ws.Range("K73").Value = 125002000
$ws.Range("L73").Value = 33336672
$ws.Range("M73").Value = -125001064
$ws.Range("N73").Value = -33338544
$ws.Range("H80").Value = 4298.3335
$ws.Range("I80").Value = 3047.5
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 3047.5
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -2049.5
$ws.Range("N80").Value = -8796
$ws.Range("H83").Value = 4298.3335
$ws.Range("I83").Value = 3047.5
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 15237.5
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -10245.5
$ws.Range("N83").Value = -43984
$ws.Range("H132").Value = 2173.6667
$ws.Range("I132").Value = 1881.1177
$ws.Range("K132").Value = 5643.3531
$ws.Range("M132").Value = -3113.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1743.1578
$ws.Range("I68").Value = 1613.875
$ws.Range("K68").Value = 1613.875
$ws.Range("M68").Value = -864.875
$ws.Range("H71").Value = 1743.1578
$ws.Range("I71").Value = 1613.875
$ws.Range("K71").Value = 8069.375
$ws.Range("M71").Value = -4325.375
$ws.Range("H82").Value = 1933
$ws.Range("J82").Value = 1949.5
$ws.Range("L82").Value = 1949.5
$ws.Range("N82").Value = -2671.5
$ws.Range("H85").Value = 1933
$ws.Range("J85").Value = 1949.5
$ws.Range("L85").Value = 1949.5
$ws.Range("N85").Value = -4445.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H96").Value = 828.4167
$ws.Range("I96").Value = 712.7646999999999
$ws.Range("K96").Value = 712.7646999999999
$ws.Range("M96").Value = 660.2353000000001
